$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 8.908234419588917
$ws.Cells.Item(2, 3).Value = 5.734832280251772
$ws.Cells.Item(2, 5).Value = 24.90221787668539
$ws.Cells.Item(2, 6).Value = 39.0560475493529
$ws.Cells.Item(2, 7).Value = 23.09524334873279
$ws.Cells.Item(2, 8).Value = 13.06200285466351
$ws.Cells.Item(2, 9).Value = 18.51494441848695
$ws.Cells.Item(2, 10).Value = 7.66539702960366
$ws.Cells.Item(2, 11).Value = 8.706436435789486
$ws.Cells.Item(2, 15).Value = 19.09416426328727

$ws.Cells.Item(3, 2).Value = 8.50461753687941
$ws.Cells.Item(3, 3).Value = 5.547308822190931
$ws.Cells.Item(3, 5).Value = 24.54041403265173
$ws.Cells.Item(3, 6).Value = 38.90334412245793
$ws.Cells.Item(3, 7).Value = 23.28705142676716
$ws.Cells.Item(3, 8).Value = 13.12217251779947
$ws.Cells.Item(3, 9).Value = 18.63393722250179
$ws.Cells.Item(3, 10).Value = 7.68836336948521
$ws.Cells.Item(3, 11).Value = 8.378345038566485
$ws.Cells.Item(3, 15).Value = 19.20928200101305

$ws.Cells.Item(4, 2).Value = 8.246790415916285
$ws.Cells.Item(4, 3).Value = 5.428286289247843
$ws.Cells.Item(4, 5).Value = 24.3212268511024
$ws.Cells.Item(4, 6).Value = 38.82061667655109
$ws.Cells.Item(4, 7).Value = 23.41345846946988
$ws.Cells.Item(4, 8).Value = 13.16120139738486
$ws.Cells.Item(4, 9).Value = 18.71078442813988
$ws.Cells.Item(4, 10).Value = 7.703234902668191
$ws.Cells.Item(4, 11).Value = 8.168877374069588
$ws.Cells.Item(4, 15).Value = 19.28425423376626

$ws.Cells.Item(5, 2).Value = 8.139333991487408
$ws.Cells.Item(5, 3).Value = 5.378860674170996
$ws.Cells.Item(5, 5).Value = 24.23275944726774
$ws.Cells.Item(5, 6).Value = 38.78970312229546
$ws.Cells.Item(5, 7).Value = 23.46713242661399
$ws.Cells.Item(5, 8).Value = 13.17763093657468
$ws.Cells.Item(5, 9).Value = 18.74305416181851
$ws.Cells.Item(5, 10).Value = 7.709489318918407
$ws.Cells.Item(5, 11).Value = 8.081586531448879
$ws.Cells.Item(5, 15).Value = 19.31588465906108

$ws.Cells.Item(6, 2).Value = 8.121350650040188
$ws.Cells.Item(6, 3).Value = 5.370599503115963
$ws.Cells.Item(6, 5).Value = 24.21812417230191
$ws.Cells.Item(6, 6).Value = 38.78473960447521
$ws.Cells.Item(6, 7).Value = 23.47617525492986
$ws.Cells.Item(6, 8).Value = 13.18039077873686
$ws.Cells.Item(6, 9).Value = 18.7484701993297
$ws.Cells.Item(6, 10).Value = 7.710539600516101
$ws.Cells.Item(6, 11).Value = 8.066978017071582
$ws.Cells.Item(6, 15).Value = 19.32120201533074

$ws.Cells.Item(7, 2).Value = 8.245350716858681
$ws.Cells.Item(7, 3).Value = 5.427623380934434
$ws.Cells.Item(7, 5).Value = 24.32003015216562
$ws.Cells.Item(7, 6).Value = 38.82018840539271
$ws.Cells.Item(7, 7).Value = 23.41417359359479
$ws.Cells.Item(7, 8).Value = 13.16142084531343
$ws.Cells.Item(7, 9).Value = 18.71121576410338
$ws.Cells.Item(7, 10).Value = 7.703318465164559
$ws.Cells.Item(7, 11).Value = 8.167707837835314
$ws.Cells.Item(7, 15).Value = 19.28467644535808

$ws.Cells.Item(8, 2).Value = 8.771213727039271
$ws.Cells.Item(8, 3).Value = 5.671010346177006
$ws.Cells.Item(8, 5).Value = 24.77691852272099
$ws.Cells.Item(8, 6).Value = 39.001121207982
$ws.Cells.Item(8, 7).Value = 23.15958053579977
$ws.Cells.Item(8, 8).Value = 13.08231720426015
$ws.Cells.Item(8, 9).Value = 18.55518877758743
$ws.Cells.Item(8, 10).Value = 7.673156286676074
$ws.Cells.Item(8, 11).Value = 8.595022619645082
$ws.Cells.Item(8, 15).Value = 19.13296657346037

$ws.Cells.Item(9, 2).Value = 9.718403317534646
$ws.Cells.Item(9, 3).Value = 6.115405494805129
$ws.Cells.Item(9, 5).Value = 25.69157923777959
$ws.Cells.Item(9, 6).Value = 39.44222561581658
$ws.Cells.Item(9, 7).Value = 22.72928500890836
$ws.Cells.Item(9, 8).Value = 12.94369786071288
$ws.Cells.Item(9, 9).Value = 18.27915844547184
$ws.Cells.Item(9, 10).Value = 7.620096133705443
$ws.Cells.Item(9, 11).Value = 9.366132423431003
$ws.Cells.Item(9, 15).Value = 18.86949883738896

$ws.Cells.Item(10, 2).Value = 10.3579522645164
$ws.Cells.Item(10, 3).Value = 6.41941267785513
$ws.Cells.Item(10, 5).Value = 26.36857484659921
$ws.Cells.Item(10, 6).Value = 39.81698714206543
$ws.Cells.Item(10, 7).Value = 22.45580070434747
$ws.Cells.Item(10, 8).Value = 12.85186323154614
$ws.Cells.Item(10, 9).Value = 18.09447546568441
$ws.Cells.Item(10, 10).Value = 7.584792328307165
$ws.Cells.Item(10, 11).Value = 9.888231751830157
$ws.Cells.Item(10, 15).Value = 18.69667786756109

$ws.Cells.Item(11, 2).Value = 10.63584765781438
$ws.Cells.Item(11, 3).Value = 6.552404954102951
$ws.Cells.Item(11, 5).Value = 26.67630885630599
$ws.Cells.Item(11, 6).Value = 39.99799540779984
$ws.Cells.Item(11, 7).Value = 22.34078541206623
$ws.Cells.Item(11, 8).Value = 12.81224802439765
$ws.Cells.Item(11, 9).Value = 18.01436398930595
$ws.Cells.Item(11, 10).Value = 7.569523950455189
$ws.Cells.Item(11, 11).Value = 10.11548644326757
$ws.Cells.Item(11, 15).Value = 18.62256417476247

$ws.Cells.Item(12, 2).Value = 10.73915092279365
$ws.Cells.Item(12, 3).Value = 6.601973228764582
$ws.Cells.Item(12, 5).Value = 26.79269787609531
$ws.Cells.Item(12, 6).Value = 40.06800637720949
$ws.Cells.Item(12, 7).Value = 22.29859502985593
$ws.Cells.Item(12, 8).Value = 12.79755677717737
$ws.Cells.Item(12, 9).Value = 17.98458690900094
$ws.Cells.Item(12, 10).Value = 7.563855529519954
$ws.Cells.Item(12, 11).Value = 10.20002594868309
$ws.Cells.Item(12, 15).Value = 18.59514729739691

$ws.Cells.Item(13, 2).Value = 10.716989313589
$ws.Cells.Item(13, 3).Value = 6.591333530709149
$ws.Cells.Item(13, 5).Value = 26.76763949213143
$ws.Cells.Item(13, 6).Value = 40.05286379000319
$ws.Cells.Item(13, 7).Value = 22.30762064208054
$ws.Cells.Item(13, 8).Value = 12.80070701366897
$ws.Cells.Item(13, 9).Value = 17.99097508278449
$ws.Cells.Item(13, 10).Value = 7.565071289534925
$ws.Cells.Item(13, 11).Value = 10.18188693531667
$ws.Cells.Item(13, 15).Value = 18.60102316266834

$ws.Cells.Item(14, 2).Value = 10.64438543414428
$ws.Cells.Item(14, 3).Value = 6.556499020941773
$ws.Cells.Item(14, 5).Value = 26.68588773518844
$ws.Cells.Item(14, 6).Value = 40.00372612853665
$ws.Cells.Item(14, 7).Value = 22.33728698605949
$ws.Cells.Item(14, 8).Value = 12.8110331534129
$ws.Cells.Item(14, 9).Value = 18.01190300923689
$ws.Cells.Item(14, 10).Value = 7.569055336169815
$ws.Cells.Item(14, 11).Value = 10.12247218386698
$ws.Cells.Item(14, 15).Value = 18.62029557267198

$ws.Cells.Item(15, 2).Value = 10.59966069303895
$ws.Cells.Item(15, 3).Value = 6.535057760521314
$ws.Cells.Item(15, 5).Value = 26.63579058787715
$ws.Cells.Item(15, 6).Value = 39.97381750902846
$ws.Cells.Item(15, 7).Value = 22.35563641299765
$ws.Cells.Item(15, 8).Value = 12.81739858975188
$ws.Cells.Item(15, 9).Value = 18.02479476972606
$ws.Cells.Item(15, 10).Value = 7.571510432707333
$ws.Cells.Item(15, 11).Value = 10.08588022066494
$ws.Cells.Item(15, 15).Value = 18.63218494187125

$ws.Cells.Item(16, 2).Value = 10.33952387999146
$ws.Cells.Item(16, 3).Value = 6.410611927537757
$ws.Cells.Item(16, 5).Value = 26.3484495545835
$ws.Cells.Item(16, 6).Value = 39.80536573603958
$ws.Cells.Item(16, 7).Value = 22.46350732313422
$ws.Cells.Item(16, 8).Value = 12.85449561206944
$ws.Cells.Item(16, 9).Value = 18.09978930843939
$ws.Cells.Item(16, 10).Value = 7.58580604245733
$ws.Cells.Item(16, 11).Value = 9.873170042116977
$ws.Cells.Item(16, 15).Value = 18.70161205701893

$ws.Cells.Item(17, 2).Value = 10.17655685859831
$ws.Cells.Item(17, 3).Value = 6.332887473277046
$ws.Cells.Item(17, 5).Value = 26.17203615969654
$ws.Cells.Item(17, 6).Value = 39.70469004505246
$ws.Cells.Item(17, 7).Value = 22.5320975149943
$ws.Cells.Item(17, 8).Value = 12.87780647228855
$ws.Cells.Item(17, 9).Value = 18.14679415706179
$ws.Cells.Item(17, 10).Value = 7.594778347378478
$ws.Cells.Item(17, 11).Value = 9.740021221514425
$ws.Cells.Item(17, 15).Value = 18.74535730829115

$ws.Cells.Item(18, 2).Value = 10.08159790185772
$ws.Cells.Item(18, 3).Value = 6.287684938409043
$ws.Cells.Item(18, 5).Value = 26.07055211010046
$ws.Cells.Item(18, 6).Value = 39.64777793551252
$ws.Cells.Item(18, 7).Value = 22.57243222432831
$ws.Cells.Item(18, 8).Value = 12.89141768114797
$ws.Cells.Item(18, 9).Value = 18.17419744974448
$ws.Cells.Item(18, 10).Value = 7.600013505228419
$ws.Cells.Item(18, 11).Value = 9.662474811627922
$ws.Cells.Item(18, 15).Value = 18.77094225188439

$ws.Cells.Item(19, 2).Value = 10.04923792277135
$ws.Cells.Item(19, 3).Value = 6.272295687287663
$ws.Cells.Item(19, 5).Value = 26.03619203985236
$ws.Cells.Item(19, 6).Value = 39.62868054869548
$ws.Cells.Item(19, 7).Value = 22.58624021918387
$ws.Cells.Item(19, 8).Value = 12.89606115718682
$ws.Cells.Item(19, 9).Value = 18.1835388775379
$ws.Cells.Item(19, 10).Value = 7.601798852857167
$ws.Cells.Item(19, 11).Value = 9.636054989353168
$ws.Cells.Item(19, 15).Value = 18.77967764542727

$ws.Cells.Item(20, 2).Value = 10.19403214787103
$ws.Cells.Item(20, 3).Value = 6.341213095339731
$ws.Cells.Item(20, 5).Value = 26.19081809786371
$ws.Cells.Item(20, 6).Value = 39.71530460990019
$ws.Cells.Item(20, 7).Value = 22.52470444877639
$ws.Cells.Item(20, 8).Value = 12.87530394243093
$ws.Cells.Item(20, 9).Value = 18.14175240385226
$ws.Cells.Item(20, 10).Value = 7.593815520269549
$ws.Cells.Item(20, 11).Value = 9.754295109946584
$ws.Cells.Item(20, 15).Value = 18.74065668119193

$ws.Cells.Item(21, 2).Value = 10.66576371512473
$ws.Cells.Item(21, 3).Value = 6.566752508202045
$ws.Cells.Item(21, 5).Value = 26.70990494889617
$ws.Cells.Item(21, 6).Value = 40.01811962236652
$ws.Cells.Item(21, 7).Value = 22.32853615743727
$ws.Cells.Item(21, 8).Value = 12.80799170215859
$ws.Cells.Item(21, 9).Value = 18.00574080113748
$ws.Cells.Item(21, 10).Value = 7.567882051367489
$ws.Cells.Item(21, 11).Value = 10.13996521595011
$ws.Cells.Item(21, 15).Value = 18.6146171954426

$ws.Cells.Item(22, 2).Value = 10.96280073567897
$ws.Cells.Item(22, 3).Value = 6.709521645621243
$ws.Cells.Item(22, 5).Value = 27.04827542716474
$ws.Cells.Item(22, 6).Value = 40.22455420751788
$ws.Cells.Item(22, 7).Value = 22.20828190363477
$ws.Cells.Item(22, 8).Value = 12.76580694323655
$ws.Cells.Item(22, 9).Value = 17.92010920781791
$ws.Cells.Item(22, 10).Value = 7.551593711993809
$ws.Cells.Item(22, 11).Value = 10.38316554756458
$ws.Cells.Item(22, 15).Value = 18.53602272684428

$ws.Cells.Item(23, 2).Value = 10.8053131543952
$ws.Cells.Item(23, 3).Value = 6.633756223803287
$ws.Cells.Item(23, 5).Value = 26.8677967672269
$ws.Cells.Item(23, 6).Value = 40.11361212038782
$ws.Cells.Item(23, 7).Value = 22.27173192554458
$ws.Cells.Item(23, 8).Value = 12.78815651595105
$ws.Cells.Item(23, 9).Value = 17.96551466396696
$ws.Cells.Item(23, 10).Value = 7.560226794923405
$ws.Cells.Item(23, 11).Value = 10.25418798791854
$ws.Cells.Item(23, 15).Value = 18.57762395660592

$ws.Cells.Item(24, 2).Value = 10.18613551190611
$ws.Cells.Item(24, 3).Value = 6.337450691827342
$ws.Cells.Item(24, 5).Value = 26.18232696719434
$ws.Cells.Item(24, 6).Value = 39.71050274726306
$ws.Cells.Item(24, 7).Value = 22.52804404745404
$ws.Cells.Item(24, 8).Value = 12.87643468359209
$ws.Cells.Item(24, 9).Value = 18.14403059796789
$ws.Cells.Item(24, 10).Value = 7.594250574956593
$ws.Cells.Item(24, 11).Value = 9.747844988317432
$ws.Cells.Item(24, 15).Value = 18.74278047906702

$ws.Cells.Item(25, 2).Value = 9.471742947419029
$ws.Cells.Item(25, 3).Value = 5.998975749987812
$ws.Cells.Item(25, 5).Value = 25.44280242791172
$ws.Cells.Item(25, 6).Value = 39.31384082374133
$ws.Cells.Item(25, 7).Value = 22.83824645152026
$ws.Cells.Item(25, 8).Value = 12.97943675077759
$ws.Cells.Item(25, 9).Value = 18.35064094092346
$ws.Cells.Item(25, 10).Value = 7.633801877723469
$ws.Cells.Item(25, 11).Value = 9.165096910469202
$ws.Cells.Item(25, 15).Value = 18.93713046528601
